# Automatische test-sync: 2025-08-05 16:57:50
# Append a new mail-log entry to the "Logs" sheet and refresh the
# "Dashboard" summary sheet accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new row (row 10) to the Logs sheet -----------------------
$newRow = 10

$logs.Cells.Item($newRow, 1).Value = "Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 16:57:36"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend the conditional formatting ranges to cover the new row -------
# The sheet has per-column conditional-formatting rules scoped to rows 2-9;
# now that row 10 holds data too, every rule's range grows by one row.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $col + "2:" + $col + "9"
    $newRange = $col + "2:" + $col + "10"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Update the Dashboard summary table -----------------------------------
# "Planning / Afspraak" moves to row 2 with an incremented count, and
# "Klantenservice / Contact" moves to row 3, keeping its previous count.
$dashboard.Cells.Item(2, 1).Value = "Planning / Afspraak"
$dashboard.Cells.Item(2, 2).Value = 4

$dashboard.Cells.Item(3, 1).Value = "Klantenservice / Contact"
$dashboard.Cells.Item(3, 2).Value = 3
